# Updates the cryptocurrency price/volume table to reflect the latest
# scrape (columns D = Price, E = Volume(1h)). Values are kept as literal
# text (matching the workbook's original inlineStr cells) so that figures
# like trailing zeros ("29.20") or percent signs ("0.78%") are preserved
# exactly instead of being reinterpreted by Excel as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "288.42" },
    @{ Cell = "E2"; Value = "0.78%" },
    @{ Cell = "D3"; Value = "29.20" },
    @{ Cell = "E3"; Value = "1.38%" },
    @{ Cell = "D4"; Value = "5.269" },
    @{ Cell = "E4"; Value = "4.67%" },
    @{ Cell = "D5"; Value = "0.06985" },
    @{ Cell = "E5"; Value = "4.15%" },
    @{ Cell = "D6"; Value = "7.447" },
    @{ Cell = "E6"; Value = "1.28%" },
    @{ Cell = "D7"; Value = "3.558" },
    @{ Cell = "E7"; Value = "5.20%" },
    @{ Cell = "D8"; Value = "1.393" },
    @{ Cell = "E8"; Value = "1.58%" },
    @{ Cell = "D9"; Value = "0.9042" },
    @{ Cell = "E9"; Value = "-4.07%" },
    @{ Cell = "D10"; Value = "0.1599" },
    @{ Cell = "E10"; Value = "2.03%" },
    @{ Cell = "D11"; Value = "0.07597" },
    @{ Cell = "E11"; Value = "13.06%" },
    @{ Cell = "D12"; Value = "0.07726" },
    @{ Cell = "E12"; Value = "2.34%" },
    @{ Cell = "D13"; Value = "0.02914" },
    @{ Cell = "E13"; Value = "-2.02%" },
    @{ Cell = "D14"; Value = "0.09023" },
    @{ Cell = "E14"; Value = "0.27%" },
    @{ Cell = "D15"; Value = "0.001575" },
    @{ Cell = "E15"; Value = "-1.71%" },
    @{ Cell = "D16"; Value = "0.0006510" },
    @{ Cell = "E16"; Value = "0.64%" },
    @{ Cell = "D17"; Value = "0.006145" },
    @{ Cell = "E17"; Value = "-6.39%" },
    @{ Cell = "D18"; Value = "3.482" },
    @{ Cell = "E18"; Value = "-0.38%" },
    @{ Cell = "D19"; Value = "2.231" },
    @{ Cell = "E19"; Value = "-0.69%" },
    @{ Cell = "D20"; Value = "0.3245" },
    @{ Cell = "E20"; Value = "1.13%" },
    @{ Cell = "E21"; Value = "2.02%" },
    @{ Cell = "D22"; Value = "4.009" },
    @{ Cell = "E22"; Value = "-2.01%" },
    @{ Cell = "D23"; Value = "0.1599" },
    @{ Cell = "E23"; Value = "3.28%" },
    @{ Cell = "E24"; Value = "0.97%" },
    @{ Cell = "D25"; Value = "0.001211" },
    @{ Cell = "E25"; Value = "2.76%" },
    @{ Cell = "D26"; Value = "0.004151" },
    @{ Cell = "E26"; Value = "-7.83%" },
    @{ Cell = "E27"; Value = "-6.11%" },
    @{ Cell = "E28"; Value = "3.41%" },
    @{ Cell = "D40"; Value = "0.04356" },
    @{ Cell = "E40"; Value = "3.51%" },
    @{ Cell = "D41"; Value = "0.006971" },
    @{ Cell = "E41"; Value = "3.35%" },
    @{ Cell = "D42"; Value = "0.1248" },
    @{ Cell = "E42"; Value = "-0.81%" },
    @{ Cell = "D43"; Value = "0.002069" },
    @{ Cell = "E43"; Value = "2.80%" },
    @{ Cell = "D44"; Value = "0.01168" },
    @{ Cell = "E44"; Value = "-4.45%" },
    @{ Cell = "D45"; Value = "0.00005872" },
    @{ Cell = "E45"; Value = "5.49%" },
    @{ Cell = "D47"; Value = "0.01299" },
    @{ Cell = "E47"; Value = "-0.29%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force Text format first so the numeric-looking / percent strings
    # are stored as-is instead of being parsed into a Double.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    # Restore the default (unstyled) look the cell had before, since
    # only the displayed text should change, not the cell formatting.
    $cell.Style = "Normal"
}
